# "moved figure around a bit"
# - Update the column-header textbox spacing/text on slide 1.
# - Nudge ten of the picture shapes in the Figure 3 grid to new positions.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Text Box 2 (column headers) ---------------------------------------
$s.Shapes.Item("Text Box 2").TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "                Pgp3 MBA                             CT694                               LFA Latex                             LFA Gold                                       LFA Field"

# --- Picture position nudges --------------------------------------------
# Values below are expressed in points (EMU / 12700) with enough decimal
# precision to land on the exact target EMU after COM's internal
# (single-precision) rounding.

# Picture 22 : (1362059,1481471) -> (1373119,1614495)
$sh = $s.Shapes.Item("Picture 22")
$sh.Left = 108.11961
$sh.Top = 127.1256

# Picture 24 : (1362059,2702481) -> (1373120,2735884)
$sh = $s.Shapes.Item("Picture 24")
$sh.Left = 108.1197
$sh.Top = 215.424

# Picture 30 : (3191277,1500900) -> (3197008,1612745)
$sh = $s.Shapes.Item("Picture 30")
$sh.Left = 251.73292
$sh.Top = 126.9878

# Picture 32 : (3185343,2684651) -> (3191277,2741816)
$sh = $s.Shapes.Item("Picture 32")
$sh.Left = 251.2817
$sh.Top = 215.8911

# Picture 38 : (4997548,1499826) -> (5020897,1612745)
$sh = $s.Shapes.Item("Picture 38")
$sh.Left = 395.34621
$sh.Top = 126.9878

# Picture 40 : (4997548,2684651) -> (5020897,2744917)
$sh = $s.Shapes.Item("Picture 40")
$sh.Left = 395.34621
$sh.Top = 216.13521

# Picture 42 : (4997548,3884016) -> (4994773,3880145)
$sh = $s.Shapes.Item("Picture 42")
$sh.Left = 393.2892
$sh.Top = 305.5233

# Picture 48 : (6809795,3876846) -> (6815686,3877295)
$sh = $s.Shapes.Item("Picture 48")
$sh.Left = 536.6682
$sh.Top = 305.2989

# Picture 50 : (8621957,3861851) -> (8633824,3877295)
$sh = $s.Shapes.Item("Picture 50")
$sh.Left = 679.8287
$sh.Top = 305.2989

# Picture 54 : (8622022,4999959) -> (8633824,4999959)
$sh = $s.Shapes.Item("Picture 54")
$sh.Left = 679.8287
$sh.Top = 393.6976
